$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 38-39 for the new "obsidian" igneous rock entry,
# pushing the existing rows 38-49 (breccia..limestoneDesc) down to 40-51.
$ws.Rows("38:39").Insert()

# New row 38: obsidian key/value, VoiceDuration 1
$ws.Cells.Item(38, 1).Value = "obsidian"
$ws.Cells.Item(38, 2).Value = "Obsidian"
$ws.Cells.Item(38, 3).Value = 1

# New row 39: obsidianDesc key/value, VoiceDuration 5
$ws.Cells.Item(39, 1).Value = "obsidianDesc"
$ws.Cells.Item(39, 2).Value = "Igneous rock. Extrusive"
$ws.Cells.Item(39, 3).Value = 5

# Update the view to match the target state (scroll position + selection)
$ws.Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C40").Select()
